$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1978.6957
$ws.Cells.Item(40, 9).Value = 2399.9333
$ws.Cells.Item(40, 10).Value = 1188.875
$ws.Cells.Item(40, 11).Value = 2399.9333
$ws.Cells.Item(40, 12).Value = 1188.875
$ws.Cells.Item(40, 13).Value = -2224.9333
$ws.Cells.Item(40, 14).Value = -1538.875
$ws.Cells.Item(64, 8).Value = 3301.9697
$ws.Cells.Item(64, 9).Value = 2829.2104
$ws.Cells.Item(64, 10).Value = 3943.5715
$ws.Cells.Item(64, 11).Value = 2829.2104
$ws.Cells.Item(64, 12).Value = 3943.5715
$ws.Cells.Item(64, 13).Value = -2581.2104
$ws.Cells.Item(64, 14).Value = -4439.5715
$ws.Cells.Item(67, 8).Value = 3301.9697
$ws.Cells.Item(67, 9).Value = 2829.2104
$ws.Cells.Item(67, 10).Value = 3943.5715
$ws.Cells.Item(67, 11).Value = 2829.2104
$ws.Cells.Item(67, 12).Value = 3943.5715
$ws.Cells.Item(67, 13).Value = -1971.2104
$ws.Cells.Item(67, 14).Value = -5659.5715
$ws.Cells.Item(74, 8).Value = 5233.1665
$ws.Cells.Item(74, 9).Value = 4879.8
$ws.Cells.Item(74, 10).Value = 7000
$ws.Cells.Item(74, 11).Value = 4879.8
$ws.Cells.Item(74, 12).Value = 7000
$ws.Cells.Item(74, 13).Value = -3943.8
$ws.Cells.Item(74, 14).Value = -8872
$ws.Cells.Item(77, 8).Value = 5233.1665
$ws.Cells.Item(77, 9).Value = 4879.8
$ws.Cells.Item(77, 10).Value = 7000
$ws.Cells.Item(77, 11).Value = 24399
$ws.Cells.Item(77, 12).Value = 35000
$ws.Cells.Item(77, 13).Value = -19719
$ws.Cells.Item(77, 14).Value = -44360
$ws.Cells.Item(100, 8).Value = 1994.3334
$ws.Cells.Item(100, 9).Value = 1741.5
$ws.Cells.Item(100, 10).Value = 2500
$ws.Cells.Item(100, 11).Value = 1741.5
$ws.Cells.Item(100, 12).Value = 2500
$ws.Cells.Item(100, 13).Value = -1200.5
$ws.Cells.Item(100, 14).Value = -3582
$ws.Cells.Item(127, 8).Value = 52632884
$ws.Cells.Item(127, 9).Value = 483.7143
$ws.Cells.Item(127, 10).Value = 83335120
$ws.Cells.Item(127, 11).Value = 1451.1429
$ws.Cells.Item(127, 12).Value = 250005360
$ws.Cells.Item(127, 13).Value = 3508.8571
$ws.Cells.Item(127, 14).Value = -250015280
$ws.Cells.Item(131, 8).Value = 4419.8184
$ws.Cells.Item(131, 9).Value = 209.44444
$ws.Cells.Item(131, 10).Value = 7334.6924
$ws.Cells.Item(131, 11).Value = 628.33332
$ws.Cells.Item(131, 12).Value = 22004.0772
$ws.Cells.Item(131, 13).Value = 4411.66668
$ws.Cells.Item(131, 14).Value = -32084.0772
$ws.Cells.Item(135, 8).Value = 35714980
$ws.Cells.Item(135, 9).Value = 13158589
$ws.Cells.Item(135, 10).Value = 250000700
$ws.Cells.Item(135, 11).Value = 118427301
$ws.Cells.Item(135, 12).Value = 2250006300
$ws.Cells.Item(135, 13).Value = -118424766
$ws.Cells.Item(135, 14).Value = -2250011370
$ws.Cells.Item(137, 8).Value = 1778.9056
$ws.Cells.Item(137, 9).Value = 1308.8918
$ws.Cells.Item(137, 10).Value = 2865.8125
$ws.Cells.Item(137, 11).Value = 3926.6754
$ws.Cells.Item(137, 12).Value = 8597.4375
$ws.Cells.Item(137, 13).Value = -1376.6754
$ws.Cells.Item(137, 14).Value = -13697.4375
$ws.Cells.Item(141, 8).Value = 2166.8965
$ws.Cells.Item(141, 9).Value = 1235.3334
$ws.Cells.Item(141, 10).Value = 5391.5386
$ws.Cells.Item(141, 11).Value = 3706.0002
$ws.Cells.Item(141, 12).Value = 16174.6158
$ws.Cells.Item(141, 13).Value = 1473.9998
$ws.Cells.Item(141, 14).Value = -26534.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18096.562
$ws.Cells.Item(32, 9).Value = 18613.062
$ws.Cells.Item(32, 10).Value = 14423.667
$ws.Cells.Item(32, 11).Value = 18613.062
$ws.Cells.Item(32, 12).Value = 14423.667
$ws.Cells.Item(32, 13).Value = -18326.062
$ws.Cells.Item(32, 14).Value = -14997.667
$ws.Cells.Item(45, 8).Value = 1332.2034
$ws.Cells.Item(45, 9).Value = 1260.102
$ws.Cells.Item(45, 10).Value = 1685.5
$ws.Cells.Item(45, 11).Value = 1260.102
$ws.Cells.Item(45, 12).Value = 1685.5
$ws.Cells.Item(45, 13).Value = -883.1020000000001
$ws.Cells.Item(45, 14).Value = -2439.5
$ws.Cells.Item(61, 8).Value = 6582.0176
$ws.Cells.Item(61, 9).Value = 3490.383
$ws.Cells.Item(61, 10).Value = 21112.7
$ws.Cells.Item(61, 11).Value = 3490.383
$ws.Cells.Item(61, 12).Value = 21112.7
$ws.Cells.Item(61, 13).Value = -3278.383
$ws.Cells.Item(61, 14).Value = -21536.7
$ws.Cells.Item(74, 8).Value = 4488.8335
$ws.Cells.Item(74, 9).Value = 2039.56
$ws.Cells.Item(74, 10).Value = 10055.363
$ws.Cells.Item(74, 11).Value = 2039.56
$ws.Cells.Item(74, 12).Value = 10055.363
$ws.Cells.Item(74, 13).Value = -1165.56
$ws.Cells.Item(74, 14).Value = -11803.363
$ws.Cells.Item(77, 8).Value = 4488.8335
$ws.Cells.Item(77, 9).Value = 2039.56
$ws.Cells.Item(77, 10).Value = 10055.363
$ws.Cells.Item(77, 11).Value = 10197.8
$ws.Cells.Item(77, 12).Value = 50276.815
$ws.Cells.Item(77, 13).Value = -5829.799999999999
$ws.Cells.Item(77, 14).Value = -59012.815
$ws.Cells.Item(136, 8).Value = 6582.0176
$ws.Cells.Item(136, 9).Value = 3490.383
$ws.Cells.Item(136, 10).Value = 21112.7
$ws.Cells.Item(136, 11).Value = 10471.149
$ws.Cells.Item(136, 12).Value = 63338.10000000001
$ws.Cells.Item(136, 13).Value = -7921.148999999999
$ws.Cells.Item(136, 14).Value = -68438.10000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(55, 8).Value = 60700
$ws.Cells.Item(55, 10).Value = 60700
$ws.Cells.Item(55, 12).Value = 60700
$ws.Cells.Item(55, 14).Value = -61246
$ws.Cells.Item(86, 8).Value = 1690.303
$ws.Cells.Item(86, 9).Value = 1700
$ws.Cells.Item(86, 10).Value = 1593.3334
$ws.Cells.Item(86, 11).Value = 1700
$ws.Cells.Item(86, 12).Value = 1593.3334
$ws.Cells.Item(86, 13).Value = -577
$ws.Cells.Item(86, 14).Value = -3839.3334
$ws.Cells.Item(89, 8).Value = 1690.303
$ws.Cells.Item(89, 9).Value = 1700
$ws.Cells.Item(89, 10).Value = 1593.3334
$ws.Cells.Item(89, 11).Value = 8500
$ws.Cells.Item(89, 12).Value = 7966.666999999999
$ws.Cells.Item(89, 13).Value = -2884
$ws.Cells.Item(89, 14).Value = -19198.667
$ws.Cells.Item(126, 8).Value = 30000
$ws.Cells.Item(126, 10).Value = 30000
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 14).Value = -39880
$ws.Cells.Item(127, 8).Value = 63450
$ws.Cells.Item(127, 10).Value = 63450
$ws.Cells.Item(127, 12).Value = 63450
$ws.Cells.Item(127, 14).Value = -73370

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2055.3606
$ws.Cells.Item(31, 9).Value = 1516.3556
$ws.Cells.Item(31, 10).Value = 3571.3125
$ws.Cells.Item(31, 11).Value = 1516.3556
$ws.Cells.Item(31, 12).Value = 3571.3125
$ws.Cells.Item(31, 13).Value = -1221.3556
$ws.Cells.Item(31, 14).Value = -4161.3125
$ws.Cells.Item(34, 8).Value = 2055.3606
$ws.Cells.Item(34, 9).Value = 1516.3556
$ws.Cells.Item(34, 10).Value = 3571.3125
$ws.Cells.Item(34, 11).Value = 1516.3556
$ws.Cells.Item(34, 12).Value = 3571.3125
$ws.Cells.Item(34, 13).Value = -1314.3556
$ws.Cells.Item(34, 14).Value = -3975.3125
$ws.Cells.Item(58, 8).Value = 1152431.8
$ws.Cells.Item(58, 9).Value = 1684730.9
$ws.Cells.Item(58, 10).Value = 2665.6
$ws.Cells.Item(58, 11).Value = 1684730.9
$ws.Cells.Item(58, 12).Value = 2665.6
$ws.Cells.Item(58, 13).Value = -1684527.9
$ws.Cells.Item(58, 14).Value = -3071.6
$ws.Cells.Item(99, 8).Value = 4055
$ws.Cells.Item(99, 9).Value = 3076.375
$ws.Cells.Item(99, 10).Value = 5173.4287
$ws.Cells.Item(99, 11).Value = 3076.375
$ws.Cells.Item(99, 12).Value = 5173.4287
$ws.Cells.Item(99, 13).Value = -1578.375
$ws.Cells.Item(99, 14).Value = -8169.4287
$ws.Cells.Item(122, 8).Value = 8656.888999999999
$ws.Cells.Item(122, 9).Value = 9130.286
$ws.Cells.Item(122, 11).Value = 27390.858
$ws.Cells.Item(122, 13).Value = -24940.858
$ws.Cells.Item(126, 8).Value = 4055
$ws.Cells.Item(126, 9).Value = 3076.375
$ws.Cells.Item(126, 10).Value = 5173.4287
$ws.Cells.Item(126, 11).Value = 9229.125
$ws.Cells.Item(126, 12).Value = 15520.2861
$ws.Cells.Item(126, 13).Value = -6759.125
$ws.Cells.Item(126, 14).Value = -20460.2861
$ws.Cells.Item(136, 8).Value = 1152431.8
$ws.Cells.Item(136, 9).Value = 1684730.9
$ws.Cells.Item(136, 10).Value = 2665.6
$ws.Cells.Item(136, 11).Value = 5054192.699999999
$ws.Cells.Item(136, 12).Value = 7996.799999999999
$ws.Cells.Item(136, 13).Value = -5051642.699999999
$ws.Cells.Item(136, 14).Value = -13096.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 142858860
$ws.Cells.Item(22, 10).Value = 2950
$ws.Cells.Item(22, 12).Value = 8850
$ws.Cells.Item(22, 14).Value = -9188
$ws.Cells.Item(27, 8).Value = 142858860
$ws.Cells.Item(27, 10).Value = 2950
$ws.Cells.Item(27, 12).Value = 8850
$ws.Cells.Item(27, 14).Value = -9054
$ws.Cells.Item(41, 8).Value = 359.66666
$ws.Cells.Item(41, 9).Value = 89.5
$ws.Cells.Item(41, 10).Value = 900
$ws.Cells.Item(41, 11).Value = 268.5
$ws.Cells.Item(41, 12).Value = 2700
$ws.Cells.Item(41, 13).Value = 69.5
$ws.Cells.Item(41, 14).Value = -3376

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 42542.5
$ws.Cells.Item(62, 10).Value = 42542.5
$ws.Cells.Item(62, 12).Value = 42542.5
$ws.Cells.Item(62, 14).Value = -43914.5
$ws.Cells.Item(65, 8).Value = 42542.5
$ws.Cells.Item(65, 10).Value = 42542.5
$ws.Cells.Item(65, 12).Value = 127627.5
$ws.Cells.Item(65, 14).Value = -134491.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 800
$ws.Cells.Item(68, 9).Value = 700
$ws.Cells.Item(68, 11).Value = 700
$ws.Cells.Item(68, 13).Value = 49
$ws.Cells.Item(71, 8).Value = 800
$ws.Cells.Item(71, 9).Value = 700
$ws.Cells.Item(71, 11).Value = 3500
$ws.Cells.Item(71, 13).Value = 244

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 8726.315000000001
$ws.Cells.Item(15, 9).Value = 8000
$ws.Cells.Item(15, 10).Value = 9971.429
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 9971.429
$ws.Cells.Item(15, 13).Value = -7712
$ws.Cells.Item(15, 14).Value = -10547.429
